# Apply crypto list refresh (prices, 1h volume %, and two swapped coin rows)
# from the scheduled GitHub Actions data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # The Price column stores plain-looking numbers (e.g. "233.27") as TEXT
    # (inlineStr) in the source data, not as numeric cells. Assigning such a
    # string straight to .Value would let Excel auto-convert it to a real
    # number (introducing float rounding noise like 233.27000000000001), so
    # we briefly force the cell to Text format, assign, then restore the
    # cell's style to Normal (removing the temporary number-format override)
    # while keeping the stored value as text.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue 'D2' '35.104.18'
$ws.Range('E2').Value = '  +1.05%  '
Set-TextValue 'D3' '1.818.21'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  +0.47%  '
Set-TextValue 'D5' '233.27'
$ws.Range('E5').Value = '  +2.33%  '
Set-TextValue 'D6' '0.617'
$ws.Range('E6').Value = '  +0.72%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  +0.46%  '
Set-TextValue 'D8' '40.70'
$ws.Range('E8').Value = '  -7.36%  '
Set-TextValue 'D9' '0.331'
$ws.Range('E9').Value = '  +9.95%  '
Set-TextValue 'D10' '0.0686'
$ws.Range('E10').Value = '  -0.18%  '
Set-TextValue 'D12' '2.080.20'
$ws.Range('E12').Value = '  -0.67%  '
Set-TextValue 'D13' '1.822.41'
$ws.Range('E13').Value = '  -0.51%  '
Set-TextValue 'D14' '11.12'
$ws.Range('E14').Value = '  -1.41%  '
Set-TextValue 'D15' '0.663'
$ws.Range('E15').Value = '  +1.32%  '
Set-TextValue 'D16' '4.66'
$ws.Range('E16').Value = '  -1.58%  '
Set-TextValue 'D17' '35.002.67'
$ws.Range('E17').Value = '  +0.83%  '
Set-TextValue 'D18' '69.59'
$ws.Range('E18').Value = '  +1.54%  '
Set-TextValue 'D19' '0.0₃0791'
$ws.Range('E19').Value = '  +0.51%  '
Set-TextValue 'D20' '239.42'
$ws.Range('E20').Value = '  -1.92%  '
Set-TextValue 'D21' '11.89'
$ws.Range('E21').Value = '  -2.76%  '
Set-TextValue 'D22' '4.69'
$ws.Range('E22').Value = '  -2.91%  '
$ws.Range('E23').Value = '  +0.35%  '
$ws.Range('E24').Value = '  +2.88%  '
Set-TextValue 'D25' '172.92'
$ws.Range('E25').Value = '  +0.67%  '
Set-TextValue 'D26' '7.90'
$ws.Range('E26').Value = '  -1.51%  '
Set-TextValue 'D27' '17.51'
$ws.Range('E27').Value = '  -1.86%  '
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('E29').Value = '  +26.52%  '
$ws.Range('E30').Value = '  +0.58%  '
Set-TextValue 'D31' '3.338.15'
$ws.Range('E31').Value = '  +37.39%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D32' '4.05'
$ws.Range('E32').Value = '  +3.11%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D33' '0.0556'
$ws.Range('E33').Value = '  +5.19%  '
Set-TextValue 'D34' '3.99'
$ws.Range('E35').Value = '  -4.60%  '
Set-TextValue 'D36' '1.15'
$ws.Range('E36').Value = '  +7.41%  '
Set-TextValue 'D37' '93.43'
$ws.Range('E37').Value = '  +3.37%  '
Set-TextValue 'D38' '0.685'
$ws.Range('E38').Value = '  +2.87%  '
Set-TextValue 'D39' '0.0194'
$ws.Range('E39').Value = '  +0.62%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D40' '1.312.82'
$ws.Range('E40').Value = '  -2.09%  '
$ws.Range('B41').Value = 'WEMIXToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D41' '1.28'
$ws.Range('E41').Value = '  +2.93%  '
Set-TextValue 'D42' '0.989'
$ws.Range('E42').Value = '  +0.28%  '
$ws.Range('E43').Value = '  -4.77%  '
Set-TextValue 'D44' '14.65'
$ws.Range('E44').Value = '  -4.26%  '
Set-TextValue 'D45' '2.47'
$ws.Range('E45').Value = '  +1.52%  '
Set-TextValue 'D46' '2.76'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('E47').Value = '  +6.16%  '
$ws.Range('E48').Value = '  -1.47%  '
Set-TextValue 'D49' '1.994.40'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('E50').Value = '  +0.35%  '
Set-TextValue 'D51' '0.0647'
$ws.Range('E51').Value = '  +4.92%  '
